# Update "TestRunner" workbook: replace the Dice-based smoke tests with the
# new CBP / PxSearch UPAX test cases (adds a JiraTicket column, drops the old
# Module/Priority columns, and appends two new PxSearch rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1): TestID | TestName | Description | Execute | JiraTicket | (blank)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "TestID"
$ws.Range("B1").Value = "TestName"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Execute"
$ws.Range("E1").Value = "JiraTicket"
$ws.Range("F1").Value = ""

# ---------------------------------------------------------------------------
# 2. Existing data rows 2-4 get new CBP content (still 6 columns wide, F blank)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "CBP001"
$ws.Range("B2").Value = "CBP Login Test"
$ws.Range("C2").Value = "Test CBP login functionality"
$ws.Range("D2").Value = "Y"
$ws.Range("E2").Value = "CBP-1234"
$ws.Range("F2").Value = ""

$ws.Range("A3").Value = "CBP002"
$ws.Range("B3").Value = "Create 1-Day Lookout"
$ws.Range("C3").Value = "Complete workflow to create 1-day lookout"
$ws.Range("D3").Value = "Y"
$ws.Range("E3").Value = "CBP-1235"
$ws.Range("F3").Value = ""

$ws.Range("A4").Value = "CBP003"
$ws.Range("B4").Value = "TECS ID Validation"
$ws.Range("C4").Value = "Validate TECS ID generation and capture"
$ws.Range("D4").Value = "Y"
$ws.Range("E4").Value = "CBP-1236"
$ws.Range("F4").Value = ""

# ---------------------------------------------------------------------------
# 3. New rows 5-6 (PxSearch UPAX cases) - copy formatting from row 4 (A:E only,
#    matching the old sheet where these rows have no F cell at all) then fill in.
# ---------------------------------------------------------------------------
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)
$ws.Range("A6:E6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A5").Value = "CBP_PXS_001"
$ws.Range("B5").Value = "PxSearch UPAX Event Creation"
$ws.Range("C5").Value = "Create UPAX event with existing event using PxSearch for Wood, Anika"
$ws.Range("D5").Value = "Y"
$ws.Range("E5").Value = "CBP-1234"

$ws.Range("A6").Value = "CBP_PXS_002"
$ws.Range("B6").Value = "PxSearch Event Verification"
$ws.Range("C6").Value = "Verify subject and delete traveler functionality"
$ws.Range("D6").Value = "Y"
$ws.Range("E6").Value = "CBP-1235"

# ---------------------------------------------------------------------------
# 4. Row heights
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 60
$ws.Rows.Item(3).RowHeight = 105
$ws.Rows.Item(5).RowHeight = 165
$ws.Rows.Item(6).RowHeight = 105

# ---------------------------------------------------------------------------
# 5. Column widths (A, B, C get explicit custom widths)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 11.833333333333334
$ws.Columns.Item(2).ColumnWidth = 25.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.666666666666666

# ---------------------------------------------------------------------------
# 6. Active cell / selection
# ---------------------------------------------------------------------------
[void]$ws.Range("O3").Select()
